$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new worksheet "ODI Batting Extra" after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "ODI Batting Extra"

# ---------------------------------------------------------------------------
# Header row - reuse the bold/bordered/centered header style already used
# by the other sheets (copy format from "Player Info"!A1) instead of
# building a brand new style.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Player Info")
$ws1.Range("A1").Copy()
$ws3.Range("A1:F1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "MATCH_CODE"
$ws3.Range("B1").Value = "BATTING_POSITION"
$ws3.Range("C1").Value = "NUM_4"
$ws3.Range("D1").Value = "NUM_6"
$ws3.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws3.Range("F1").Value = "MAN_OF_MATCH"

# ---------------------------------------------------------------------------
# Data rows. Most of these columns are text (even when they look numeric),
# so force a Text number format before writing them - only BATTING_POSITION
# (column B, row 3) holds a genuine number. (Multi-area ranges only apply
# NumberFormat to their first area in this engine, so set it per
# contiguous block.)
# ---------------------------------------------------------------------------
$ws3.Range("A2:A3").NumberFormat = "@"
$ws3.Range("B2:B2").NumberFormat = "@"
$ws3.Range("C2:E3").NumberFormat = "@"
$ws3.Range("F2:F3").NumberFormat = "@"

$ws3.Range("A2").Value = "4433"
$ws3.Range("B2").Value = "'"
$ws3.Range("C2").Value = "'"
$ws3.Range("D2").Value = "'"
$ws3.Range("E2").Value = "'"
$ws3.Range("F2").Value = "NO"

$ws3.Range("A3").Value = "4434"
$ws3.Range("B3").Value = 4
$ws3.Range("C3").Value = "3"
$ws3.Range("D3").Value = "0"
$ws3.Range("E3").Value = "4.68%"
$ws3.Range("F3").Value = "NO"
